# Regenerate save_data: replace column G ("K") values with new strikeout
# counts computed from K (strikeouts) instead of the old Strike# column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 2
    14 = 2
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 0
    24 = 2
    25 = 3
    26 = 0
    27 = 1
    28 = 1
    29 = 3
    30 = 1
    31 = 1
    32 = 4
    33 = 0
    34 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
